$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 274, shifting existing rows 274.. down by one.
$ws.Rows.Item(274).Insert()

# Copy the number format (date style) from the row above's Date cell,
# which previously lived at D274 and is now, after the insert, on D275.
$ws.Range("D275").Copy()
$ws.Range("D274").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's values.
$ws.Cells.Item(274, 1).Value = 8
$ws.Cells.Item(274, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(274, 3).Value = "Coquimbo"
$ws.Cells.Item(274, 4).Value = 44769
$ws.Cells.Item(274, 5).Value = 4
$ws.Cells.Item(274, 6).Value = 100112012
$ws.Cells.Item(274, 7).Value = "Espinaca"
$ws.Cells.Item(274, 8).Value = "Sin especificar"
$ws.Cells.Item(274, 9).Value = "Primera"
$ws.Cells.Item(274, 10).Value = 2000
$ws.Cells.Item(274, 11).Value = 500
$ws.Cells.Item(274, 12).Value = 600
$ws.Cells.Item(274, 13).Value = 550
$ws.Cells.Item(274, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(274, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(274, 16).Value = 1100
$ws.Cells.Item(274, 17).Value = 0.5
$ws.Cells.Item(274, 18).Value = "Hortaliza"
